# Daily attendance processing - 2025-11-30 13:51:45
# Normalizes the "Recorded By" (column G) entries on the
# "Session Analysis Results" sheet: each comma-separated list of
# recorders is re-ordered (case-insensitive ascending, ties broken by
# reverse of original order) so that automated "System"/"system"
# entries consistently sort after real user identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($current -eq $null) { continue }
    if ($current -eq "") { continue }

    $parts = $current -split ", "
    if ($parts.Count -le 1) { continue }

    $reversedParts = $parts[($parts.Count - 1)..0]
    $newParts = $reversedParts | Sort-Object { $_.ToLower() }
    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $current) {
        $cell.Value = $newValue
    }
}
